$wb = $excel.ActiveWorkbook

# ALC row 21
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = ""

# ALC row 23
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = ""
$ws.Range("N23").Value = ""

# ALC row 55
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 281
$ws.Range("I55").Value = 257.33334
$ws.Range("J55").Value = 311.42856
$ws.Range("K55").Value = 257.33334
$ws.Range("L55").Value = 311.42856
$ws.Range("M55").Value = -43.33334000000002
$ws.Range("N55").Value = -739.4285600000001

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7999.4
$ws.Range("I64").Value = 4998.5
$ws.Range("J64").Value = 10000
$ws.Range("K64").Value = 4998.5
$ws.Range("L64").Value = 10000
$ws.Range("M64").Value = -4750.5

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 7999.4
$ws.Range("I67").Value = 4998.5
$ws.Range("J67").Value = 10000
$ws.Range("K67").Value = 4998.5
$ws.Range("L67").Value = 10000
$ws.Range("M67").Value = -4140.5

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6805.8887
$ws.Range("I74").Value = 10938.25
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 10938.25
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -10002.25
$ws.Range("N74").Value = -5372

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5759.4
$ws.Range("I76").Value = 4776.5
$ws.Range("J76").Value = 6414.6665
$ws.Range("K76").Value = 4776.5
$ws.Range("L76").Value = 6414.6665
$ws.Range("M76").Value = -4461.5
$ws.Range("N76").Value = -7044.6665

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 6805.8887
$ws.Range("I77").Value = 10938.25
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 54691.25
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -50011.25
$ws.Range("N77").Value = -26860

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5759.4
$ws.Range("I79").Value = 4776.5
$ws.Range("J79").Value = 6414.6665
$ws.Range("K79").Value = 4776.5
$ws.Range("L79").Value = 6414.6665
$ws.Range("M79").Value = -3684.5
$ws.Range("N79").Value = -8598.666499999999

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 49989.81
$ws.Range("I112").Value = 2874.5
$ws.Range("J112").Value = 61075.766
$ws.Range("K112").Value = 8623.5
$ws.Range("L112").Value = 183227.298
$ws.Range("M112").Value = -7515.5
$ws.Range("N112").Value = -185443.298

# ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 125000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 125000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 125000
$ws.Range("N133").Value = -135120

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7940523.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 7940523.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 23821570.5
$ws.Range("M137").Value = ""
$ws.Range("N137").Value = -23826670.5

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3253.7693
$ws.Range("I63").Value = 1912.375
$ws.Range("J63").Value = 5400
$ws.Range("K63").Value = 1912.375
$ws.Range("L63").Value = 5400
$ws.Range("M63").Value = -1226.375
$ws.Range("N63").Value = -6772

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3253.7693
$ws.Range("I66").Value = 1912.375
$ws.Range("J66").Value = 5400
$ws.Range("K66").Value = 9561.875
$ws.Range("L66").Value = 27000
$ws.Range("M66").Value = -6129.875
$ws.Range("N66").Value = -33864

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1385.2
$ws.Range("I102").Value = 1273.8948
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 1273.8948
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = 348.1052

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4399.171
$ws.Range("I122").Value = 3792
$ws.Range("J122").Value = 5347.875
$ws.Range("K122").Value = 11376
$ws.Range("L122").Value = 16043.625
$ws.Range("M122").Value = -8926

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3764.9333
$ws.Range("I86").Value = 1783.1666
$ws.Range("J86").Value = 5086.1113
$ws.Range("K86").Value = 1783.1666
$ws.Range("L86").Value = 5086.1113
$ws.Range("M86").Value = -660.1666
$ws.Range("N86").Value = -7332.1113

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3764.9333
$ws.Range("I89").Value = 1783.1666
$ws.Range("J89").Value = 5086.1113
$ws.Range("K89").Value = 8915.833000000001
$ws.Range("L89").Value = 25430.5565
$ws.Range("M89").Value = -3299.833000000001
$ws.Range("N89").Value = -36662.5565

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1922
$ws.Range("I99").Value = 1277.5
$ws.Range("J99").Value = 4500
$ws.Range("K99").Value = 1277.5
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = 220.5
$ws.Range("N99").Value = -7496

# CRP row 50
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = ""

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2591.2666
$ws.Range("I105").Value = 2334.8
$ws.Range("J105").Value = 3104.2
$ws.Range("K105").Value = 2334.8
$ws.Range("L105").Value = 3104.2
$ws.Range("M105").Value = -587.8000000000002

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4644.8125
$ws.Range("I122").Value = 4286.2856
$ws.Range("J122").Value = 4923.6665
$ws.Range("K122").Value = 12858.8568
$ws.Range("L122").Value = 14770.9995
$ws.Range("M122").Value = -10408.8568
$ws.Range("N122").Value = -19670.9995

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3342
$ws.Range("I132").Value = 1012
$ws.Range("J132").Value = 4507
$ws.Range("K132").Value = 3036
$ws.Range("L132").Value = 13521
$ws.Range("M132").Value = -506

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 502989.2
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 502989.2
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 502989.2
$ws.Range("N141").Value = -513349.2

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8592.666999999999
$ws.Range("I3").Value = 8464.727999999999
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 25394.184
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = -25282.184

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1046
$ws.Range("I92").Value = 918.5
$ws.Range("J92").Value = 1155.2858
$ws.Range("K92").Value = 2755.5
$ws.Range("L92").Value = 3465.8574
$ws.Range("M92").Value = -1507.5
$ws.Range("N92").Value = -5961.857400000001

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 501532.25
$ws.Range("I117").Value = 3000
$ws.Range("J117").Value = 667709.7
$ws.Range("K117").Value = 9000
$ws.Range("L117").Value = 2003129.1
$ws.Range("M117").Value = -5558
$ws.Range("N117").Value = -2010013.1

# CUL row 130
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 1000
$ws.Range("I130").Value = 1000
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 3000
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = 2020
$ws.Range("N130").Value = ""

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 4231.316
$ws.Range("I133").Value = 4026.3333
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 12078.9999
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -7018.999899999999

# GSM row 31
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1949.5
$ws.Range("I31").Value = 1739.4
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1739.4
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1447.4

# GSM row 37
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 1949.5
$ws.Range("I37").Value = 1739.4
$ws.Range("J37").Value = 3000
$ws.Range("K37").Value = 1739.4
$ws.Range("L37").Value = 3000
$ws.Range("M37").Value = -1462.4

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4244.6787
$ws.Range("I70").Value = 3330.8333
$ws.Range("J70").Value = 4493.909
$ws.Range("K70").Value = 3330.8333
$ws.Range("L70").Value = 4493.909
$ws.Range("M70").Value = -3060.8333
$ws.Range("N70").Value = -5033.909

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4244.6787
$ws.Range("I73").Value = 3330.8333
$ws.Range("J73").Value = 4493.909
$ws.Range("K73").Value = 3330.8333
$ws.Range("L73").Value = 4493.909
$ws.Range("M73").Value = -2394.8333
$ws.Range("N73").Value = -6365.909

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1417
$ws.Range("I102").Value = 1278.3334
$ws.Range("J102").Value = 2082.6
$ws.Range("K102").Value = 1278.3334
$ws.Range("L102").Value = 2082.6
$ws.Range("M102").Value = 343.6666
$ws.Range("N102").Value = -5326.6

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2886.1
$ws.Range("I126").Value = 2566.6
$ws.Range("J126").Value = 3205.6
$ws.Range("K126").Value = 7699.799999999999
$ws.Range("L126").Value = 9616.799999999999
$ws.Range("M126").Value = -5229.799999999999
$ws.Range("N126").Value = -14556.8

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2649
$ws.Range("I132").Value = 2511.0605
$ws.Range("J132").Value = 4166.3335
$ws.Range("K132").Value = 7533.181500000001
$ws.Range("L132").Value = 12499.0005
$ws.Range("M132").Value = -5003.181500000001

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2102.389
$ws.Range("I22").Value = 1674.3636
$ws.Range("J22").Value = 2775
$ws.Range("K22").Value = 1674.3636
$ws.Range("L22").Value = 2775
$ws.Range("M22").Value = -1379.3636
$ws.Range("N22").Value = -3365

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2102.389
$ws.Range("I27").Value = 1674.3636
$ws.Range("J27").Value = 2775
$ws.Range("K27").Value = 1674.3636
$ws.Range("L27").Value = 2775
$ws.Range("M27").Value = -1567.3636
$ws.Range("N27").Value = -2989

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2035.1875
$ws.Range("I55").Value = 2482.4285
$ws.Range("J55").Value = 1687.3334
$ws.Range("K55").Value = 2482.4285
$ws.Range("L55").Value = 1687.3334
$ws.Range("M55").Value = -2309.4285

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10879.042
$ws.Range("I122").Value = 10705.389
$ws.Range("J122").Value = 11400
$ws.Range("K122").Value = 32116.167
$ws.Range("L122").Value = 34200
$ws.Range("M122").Value = -29666.167

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 15390.6
$ws.Range("I136").Value = 12988.25
$ws.Range("J136").Value = 25000
$ws.Range("K136").Value = 38964.75
$ws.Range("L136").Value = 75000
$ws.Range("M136").Value = -36414.75

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1721.3636
$ws.Range("I2").Value = 1721.3636
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1721.3636
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1609.3636
$ws.Range("N2").Value = ""

# WVR row 4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 934227.25
$ws.Range("I4").Value = 934227.25
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 934227.25
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -934114.25
$ws.Range("N4").Value = ""

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 8628
$ws.Range("I96").Value = 10623.75
$ws.Range("J96").Value = 8013.923
$ws.Range("K96").Value = 10623.75
$ws.Range("L96").Value = 8013.923
$ws.Range("M96").Value = -9250.75
$ws.Range("N96").Value = -10759.923

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5043.758
$ws.Range("I122").Value = 5428.4346
$ws.Range("J122").Value = 4159
$ws.Range("K122").Value = 16285.3038
$ws.Range("L122").Value = 12477
$ws.Range("M122").Value = -13835.3038

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1518.2667
$ws.Range("I126").Value = 1529.4615
$ws.Range("J126").Value = 1445.5
$ws.Range("K126").Value = 4588.3845
$ws.Range("L126").Value = 4336.5
$ws.Range("M126").Value = -2118.3845

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3786
$ws.Range("I132").Value = 3800.4
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 11401.2
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = -8871.200000000001
